$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.69%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.38%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.025"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.56%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07826"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.34%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.192"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-7.67%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.038"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.03%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.036"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.89%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9139"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.96%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09732"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.94%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1891"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.38%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08571"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.87%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03523"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.07%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09968"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.76%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001483"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.11%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005662"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.48%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.463"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.30%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.070"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.97%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.49%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1302"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.09%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.752"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "10.42%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2204"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.06%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04640"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.43%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.78%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004801"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.34%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.63%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "28.44%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01764"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.30%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04727"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.71%"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.29%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.04%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007659"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.55%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002180"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.19%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01039"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.70%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006041"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.50%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.437"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "135.58%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
